$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Categories")

# Rename cohort "recetox" -> "elspac" (row 39, column D label)
$ws.Cells.Item(39, 4).Value = "elspac"

# Insert a new cohort row for "genrnext" right after "genxxi" (old row 40),
# shifting urb_area_id / ath_* rows down by one.
$ws.Rows.Item(41).Insert()
$ws.Cells.Item(41, 1).Value = "cohort_id"
$ws.Cells.Item(41, 2).Value = 132
$ws.Cells.Item(41, 3).Value = $false
$ws.Cells.Item(41, 4).Value = "genrnext"

# Fix the "Rennes" urb_area_id value (was a duplicate of Grenoble's 1803);
# after the insert above it now lives on row 63.
$ws.Cells.Item(63, 2).Value = 1804

# Make "Categories" the active sheet / tab, with D42 selected.
$ws.Activate()
[void]$ws.Range("D42").Select()
